# Regenerate handoff/status report (Generate Report for Archive)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# 1) Status moved from "Ready for handoff" to "In Translation" for Test`1.md
#    and Test`2.md (both locales share the same status text in this report).
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"

$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C4").Value = "In Translation"

$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C4").Value = "In Translation"

# 2) New handoff generated for Test`1 -> "Lastest Handoff Name" column (I)
#    is refreshed with the newly generated handoff package name across all
#    rows/sheets that previously referenced the older handoff names.
$newHandoffName = "LocaleLowerCaseTest_HT_OL#Test1#20171104T090037"

$wsZhCn.Range("I3").Value = $newHandoffName
$wsZhCn.Range("I4").Value = $newHandoffName
$wsDeDe.Range("I3").Value = $newHandoffName
$wsDeDe.Range("I4").Value = $newHandoffName
